# Fruta / hortaliza, semanal
# Insert a new weekly record as row 325 in the data table, pushing all
# subsequent rows down by one (last row becomes row 399).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 325 (shifts 325..398 -> 326..399)
$ws.Rows.Item(325).Insert()

# Populate the newly inserted row with the new weekly data point
$ws.Range("A325").Value = 3
$ws.Range("B325").Value = "Femacal de La Calera"
$ws.Range("C325").Value = "Coquimbo"
$ws.Range("D325").Value = 44508
$ws.Range("E325").Value = 5
$ws.Range("F325").Value = 100112024
$ws.Range("G325").Value = "Choclo"
$ws.Range("H325").Value = "Dulce o Americano"
$ws.Range("I325").Value = "Primera"
$ws.Range("J325").Value = 85
$ws.Range("K325").Value = 40000
$ws.Range("L325").Value = 41000
$ws.Range("M325").Value = 40529
$ws.Range("N325").Value = "$/malla 70 unidades"
$ws.Range("O325").Value = "Región de Arica y Parinacota"
$ws.Range("P325").Value = 579
$ws.Range("Q325").Value = 70
$ws.Range("R325").Value = "Hortaliza"
